$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (e.g. "230.88") are stored as text, matching the source data which
# uses inline/shared strings rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.974.36'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '2.267.18'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '230.88'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").Value = '0.629'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '63.90'
$ws.Range("E7").Value = '  +3.67%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +5.87%  '
$ws.Range("D10").Value = '0.0996'
$ws.Range("E10").Value = '  +4.35%  '
$ws.Range("D11").Value = '57.40'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = '27.34'
$ws.Range("E12").Value = '  +14.97%  '
$ws.Range("E13").Value = '  +1.87%  '
$ws.Range("D14").Value = '2.605.26'
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("D15").Value = '15.73'
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").Value = '6.11'
$ws.Range("E16").Value = '  +5.07%  '
$ws.Range("D17").Value = '0.838'
$ws.Range("E17").Value = '  +3.04%  '
$ws.Range("D18").Value = '2.285.19'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").Value = '43.866.38'
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("E20").Value = '  +7.32%  '
$ws.Range("D21").Value = '73.78'
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("D22").Value = '6.12'
$ws.Range("E22").Value = '  -2.12%  '
$ws.Range("D23").Value = '252.68'
$ws.Range("E23").Value = '  +0.73%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").Value = '2.46'
$ws.Range("E25").Value = '  -4.26%  '
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("E27").Value = '  +25.23%  '
$ws.Range("D28").Value = '10.11'
$ws.Range("E28").Value = '  +2.36%  '
$ws.Range("D29").Value = '171.33'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  -1.48%  '
$ws.Range("D31").Value = '20.92'
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("E32").Value = '  -4.38%  '
$ws.Range("E33").Value = '  +2.58%  '
$ws.Range("E34").Value = '  +6.62%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").Value = '4.89'
$ws.Range("E36").Value = '  -3.26%  '
$ws.Range("D37").Value = '3.81'
$ws.Range("E37").Value = '  +5.03%  '
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("E39").Value = '  -3.63%  '
$ws.Range("E40").Value = '  +3.25%  '
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").Value = '0.000226'
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("D43").Value = '0.0992'
$ws.Range("E43").Value = '  +1.53%  '
$ws.Range("D44").Value = '17.51'
$ws.Range("E44").Value = '  +4.48%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '8.26'
$ws.Range("E45").Value = '  -5.78%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").Value = '10.48'
$ws.Range("E46").Value = '  +9.94%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").Value = '1.21'
$ws.Range("E47").Value = '  -1.02%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '98.03'
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").Value = '4.39'
$ws.Range("E49").Value = '  -3.37%  '
$ws.Range("D50").Value = '1.447.37'
$ws.Range("E50").Value = '  -1.62%  '
$ws.Range("D51").Value = '2.33'
$ws.Range("E51").Value = '  +2.95%  '

# Restore the default "Normal" style on column D so no stray
# number-format style index is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"

